$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Sponza": add a new results column K ("v1419") mirroring the existing
# J ("v1410") column - header, 10 data rows, and the 5 summary rows.
# ---------------------------------------------------------------------------
$wsSponza = $wb.Worksheets.Item("Sponza")

# Copy the formatting of column J onto the new column K first so the new
# cells pick up the right styles (header style, data style, summary style).
$wsSponza.Range("J1:J16").Copy()
$wsSponza.Range("K1:K16").PasteSpecial(-4122) # xlPasteFormats

$wsSponza.Range("K1").Value = "v1419"

$wsSponza.Range("K2").Value = 8415
$wsSponza.Range("K3").Value = 8327
$wsSponza.Range("K4").Value = 8288
$wsSponza.Range("K5").Value = 8295
$wsSponza.Range("K6").Value = 8285
$wsSponza.Range("K7").Value = 8292
$wsSponza.Range("K8").Value = 8333
$wsSponza.Range("K9").Value = 8308
$wsSponza.Range("K10").Value = 8271
$wsSponza.Range("K11").Value = 8272

$wsSponza.Range("K12").Formula = "=AVERAGE(K2:K11)"
$wsSponza.Range("K13").Formula = "=VAR.S(K2:K11)"
$wsSponza.Range("K14").Formula = "=1-T.TEST(J2:J11,K2:K11,2,3)"
$wsSponza.Range("K15").Formula = "=J12/K12"
$wsSponza.Range("K16").Formula = "=B12/K12"

# Widen the conditional formatting that highlights the DIFF ACCEPT / Perf
# ratio rows so it covers the new column too (ModifyAppliesToRange keeps the
# existing dxfId-backed rules instead of recreating them).
$wsSponza.Range("B15:J16").FormatConditions.Item(1).ModifyAppliesToRange($wsSponza.Range("B15:K16"))

# ---------------------------------------------------------------------------
# Sheet "ComplexMesh": same new column, this time labelled J (mirrors I).
# ---------------------------------------------------------------------------
$wsMesh = $wb.Worksheets.Item("ComplexMesh")

$wsMesh.Range("I1:I16").Copy()
$wsMesh.Range("J1:J16").PasteSpecial(-4122) # xlPasteFormats

$wsMesh.Range("J1").Value = "v1419"

$wsMesh.Range("J2").Value = 6422
$wsMesh.Range("J3").Value = 6327
$wsMesh.Range("J4").Value = 6296
$wsMesh.Range("J5").Value = 6329
$wsMesh.Range("J6").Value = 6281
$wsMesh.Range("J7").Value = 6304
$wsMesh.Range("J8").Value = 6366
$wsMesh.Range("J9").Value = 6346
$wsMesh.Range("J10").Value = 6360
$wsMesh.Range("J11").Value = 6326
$wsMesh.Range("J12").Value = 6277

$wsMesh.Range("J13").Formula = "=VAR.S(J2:J11)"
$wsMesh.Range("J14").Formula = "=1-T.TEST(I2:I11,J2:J11,2,3)"
$wsMesh.Range("J15").Formula = "=I12/J12"
$wsMesh.Range("J16").Formula = "=B12/J12"

$wsMesh.Range("B15:I16").FormatConditions.Item(1).ModifyAppliesToRange($wsMesh.Range("B15:J16"))

# ---------------------------------------------------------------------------
# Selection / active-tab bookkeeping: the author moved off Sponza (leaving a
# cursor at I22) and onto ComplexMesh (cursor at J20), which becomes the
# workbook's active sheet.
# ---------------------------------------------------------------------------
$wsSponza.Activate()
$wsSponza.Range("I22").Select()

$wsMesh.Activate()
$wsMesh.Range("J20").Select()
